$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 55, shifting the existing weekly records
# (old rows 55-65) down to rows 57-67.
$ws.Rows("55:56").Insert()

# New row 55: Banquete, week of 2021-11-11 (serial 44511), Volumen 300
$ws.Cells.Item(55, 1).Value = 12
$ws.Cells.Item(55, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(55, 3).Value = "Metropolitana"
$ws.Cells.Item(55, 4).Value = 44511
$ws.Cells.Item(55, 5).Value = 13
$ws.Cells.Item(55, 6).Value = 300000000
$ws.Cells.Item(55, 7).Value = "Espárragos"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Banquete"
$ws.Cells.Item(55, 10).Value = 300
$ws.Cells.Item(55, 11).Value = 1200
$ws.Cells.Item(55, 12).Value = 1200
$ws.Cells.Item(55, 13).Value = 1200
$ws.Cells.Item(55, 14).Value = "$/kilo"
$ws.Cells.Item(55, 15).Value = "Región Metropolitana"
$ws.Cells.Item(55, 16).Value = 1200
$ws.Cells.Item(55, 17).Value = 1
$ws.Cells.Item(55, 18).Value = "Hortaliza"

# New row 56: Primera, week of 2021-11-11 (serial 44511), Volumen 280
$ws.Cells.Item(56, 1).Value = 12
$ws.Cells.Item(56, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(56, 3).Value = "Metropolitana"
$ws.Cells.Item(56, 4).Value = 44511
$ws.Cells.Item(56, 5).Value = 13
$ws.Cells.Item(56, 6).Value = 300000000
$ws.Cells.Item(56, 7).Value = "Espárragos"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 280
$ws.Cells.Item(56, 11).Value = 1000
$ws.Cells.Item(56, 12).Value = 1000
$ws.Cells.Item(56, 13).Value = 1000
$ws.Cells.Item(56, 14).Value = "$/kilo"
$ws.Cells.Item(56, 15).Value = "Región Metropolitana"
$ws.Cells.Item(56, 16).Value = 1000
$ws.Cells.Item(56, 17).Value = 1
$ws.Cells.Item(56, 18).Value = "Hortaliza"
